# Apply "想去人数" (want-to-go count) updates to the 合肥-漫展信息 workbook.
# Sheet "展览"   (rows: F3,F4,F5,F6,F7,F8,F12,F18,F20,F21,F22)
# Sheet "演出"   (row : F2)
# Sheet "全部类型" (rows: F2,F4,F5,F6,F7,F8,F9,F13,F19,F21,F22,F23)

$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, [string]$addr, [double]$value)
    $ws.Range($addr).Value = $value
}

# --- 展览 (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
Set-CellValue $wsExhibit "F3"  738
Set-CellValue $wsExhibit "F4"  1454
Set-CellValue $wsExhibit "F5"  222
Set-CellValue $wsExhibit "F6"  89
Set-CellValue $wsExhibit "F7"  134
Set-CellValue $wsExhibit "F8"  6156
Set-CellValue $wsExhibit "F12" 5026
Set-CellValue $wsExhibit "F18" 59
Set-CellValue $wsExhibit "F20" 287
Set-CellValue $wsExhibit "F21" 24
Set-CellValue $wsExhibit "F22" 3530

# --- 演出 (Performance) ---
$wsShow = $wb.Worksheets.Item("演出")
Set-CellValue $wsShow "F2" 70

# --- 全部类型 (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
Set-CellValue $wsAll "F2"  70
Set-CellValue $wsAll "F4"  738
Set-CellValue $wsAll "F5"  1454
Set-CellValue $wsAll "F6"  222
Set-CellValue $wsAll "F7"  89
Set-CellValue $wsAll "F8"  134
Set-CellValue $wsAll "F9"  6156
Set-CellValue $wsAll "F13" 5026
Set-CellValue $wsAll "F19" 59
Set-CellValue $wsAll "F21" 287
Set-CellValue $wsAll "F22" 24
Set-CellValue $wsAll "F23" 3530

$wb.Save()
